# Update the classification-report table (rows 2-11) with the new model
# results. The rating-class rows have been re-ordered (A, AA, AAA, B, BB,
# BBB, C, CC, CCC, D) and the Precision / Recall / F1-Score / Support
# figures recomputed, now reported to 4 decimal places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking Precision/Recall/F1-Score/Support columns to
# stay as text (matching the original inlineStr cell type) instead of
# being auto-converted to numbers by Excel.
$ws.Range("B2:E11").NumberFormat = "@"

$ws.Range("A2").Value = "A"
$ws.Range("B2").Value = "0.6793"
$ws.Range("C2").Value = "0.5841"
$ws.Range("D2").Value = "0.6281"
$ws.Range("E2").Value = "214"

$ws.Range("A3").Value = "AA"
$ws.Range("B3").Value = "0.5882"
$ws.Range("C3").Value = "0.4762"
$ws.Range("D3").Value = "0.5263"
$ws.Range("E3").Value = "42"

$ws.Range("A4").Value = "AAA"
$ws.Range("B4").Value = "0.7917"
$ws.Range("C4").Value = "0.7037"
$ws.Range("D4").Value = "0.7451"
$ws.Range("E4").Value = "27"

$ws.Range("A5").Value = "B"
$ws.Range("B5").Value = "0.5822"
$ws.Range("C5").Value = "0.6159"
$ws.Range("D5").Value = "0.5986"
$ws.Range("E5").Value = "138"

$ws.Range("A6").Value = "BB"
$ws.Range("B6").Value = "0.6399"
$ws.Range("C6").Value = "0.6678"
$ws.Range("D6").Value = "0.6535"
$ws.Range("E6").Value = "298"

$ws.Range("A7").Value = "BBB"
$ws.Range("B7").Value = "0.6759"
$ws.Range("C7").Value = "0.7416"
$ws.Range("D7").Value = "0.7072"
$ws.Range("E7").Value = "329"

$ws.Range("A8").Value = "C"
$ws.Range("B8").Value = "0.7778"
$ws.Range("C8").Value = "1.0000"
$ws.Range("D8").Value = "0.8750"
$ws.Range("E8").Value = "7"

$ws.Range("A9").Value = "CC"
$ws.Range("B9").Value = "0.0000"
$ws.Range("C9").Value = "0.0000"
$ws.Range("D9").Value = "0.0000"
$ws.Range("E9").Value = "0"

$ws.Range("A10").Value = "CCC"
$ws.Range("B10").Value = "0.8095"
$ws.Range("C10").Value = "0.5000"
$ws.Range("D10").Value = "0.6182"
$ws.Range("E10").Value = "34"

$ws.Range("A11").Value = "D"
$ws.Range("B11").Value = "0.0000"
$ws.Range("C11").Value = "0.0000"
$ws.Range("D11").Value = "0.0000"
$ws.Range("E11").Value = "3"
